$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "Corevolume"

$values = @(
    12.930893025841099,
    6.68304685801668,
    5.12599848576673,
    4.4740996614771698,
    4.1304418589971599,
    3.9235999301504698,
    3.78784223347576,
    3.6931188904991301,
    3.6239458311972101,
    3.5716175092284801,
    11.0765504046645,
    5.6610162370444801,
    4.3113913528753303,
    3.7463357702087099,
    3.4484586264204902,
    3.2691713163816498,
    3.1514986936712099,
    3.06939399925621,
    3.0094358863216302,
    2.9640785250896502,
    9.6203896927708303,
    4.8674920789657996,
    3.6830052503435802,
    3.1870890124623599,
    2.9256596645721298,
    2.7683096780899001,
    2.6650352991681099,
    2.5929768107367299,
    2.5403550816376899,
    2.5005475783895501,
    8.4530623854175495,
    4.2378903970428796,
    3.18741211909613,
    2.7476020486192301,
    2.5157498704942101,
    2.3762018917065699,
    2.2846116020517999,
    2.2207055500723998,
    2.17403725285644,
    2.1387334263798601,
    7.5008459963847702,
    3.72913615950903,
    2.7891746904910599,
    2.3956353583476102,
    2.1881754626785801,
    2.06330879394927,
    1.98135436898015,
    1.9241716287381401,
    1.88241312489414,
    1.85082347921223,
    6.7124583287213397,
    3.3115616392606699,
    2.46401185773989,
    2.1091630836444102,
    1.9220994698217,
    1.8095089846911501,
    1.7356118545781301,
    1.6840509974618201,
    1.6463979531869399,
    1.6179140225381701
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 13).Value = $values[$i]
}

$ws.Range("M1").ColumnWidth = 10
$ws.Range("M1").Select()
